# Refactored Parser structure. Fixed problems with reading size of classes and methods.
#
# The parser that produces these workbooks previously emitted "0" for every
# class/method whose real line-count it failed to read. This edit patches in
# the now-correctly-computed line counts for the classes/methods that were
# affected (the Customer value class, its constructors, CustomerRepository,
# and the three no-arg "app/test" entry-point types).
#
# Numeric-looking values are stored as TEXT in this workbook (shared strings,
# not numbers) - e.g. "0", "3", "16" - matching the column's existing data.
# A leading apostrophe forces Excel to keep the entry as text instead of
# converting it to a number.

$wb = $excel.ActiveWorkbook

# --- classNumberOfLines: CustomerRepository class had "0" lines -> "1"
$wsClass = $wb.Worksheets.Item("classNumberOfLines")
$wsClass.Cells.Item(8, 2).Value = "'1"

# --- methodNumberOfLines: per-method line counts that were misread as "0"
$wsMethod = $wb.Worksheets.Item("methodNumberOfLines")

# PaymentAppTest()              row 6  : 0 -> 1
$wsMethod.Cells.Item(6, 3).Value = "'1"
# PaymentComponentTests()       row 10 : 0 -> 1
$wsMethod.Cells.Item(10, 3).Value = "'1"
# PaymentApp()                  row 14 : 0 -> 1
$wsMethod.Cells.Item(14, 3).Value = "'1"

# Customer getters/setters/toString, rows 15-23 : 0 -> 3
$wsMethod.Cells.Item(15, 3).Value = "'3"
$wsMethod.Cells.Item(16, 3).Value = "'3"
$wsMethod.Cells.Item(17, 3).Value = "'3"
$wsMethod.Cells.Item(18, 3).Value = "'3"
$wsMethod.Cells.Item(19, 3).Value = "'3"
$wsMethod.Cells.Item(20, 3).Value = "'3"
$wsMethod.Cells.Item(21, 3).Value = "'3"
$wsMethod.Cells.Item(22, 3).Value = "'3"
$wsMethod.Cells.Item(23, 3).Value = "'3"

# Customer()                                          row 24 : 0 -> 2
$wsMethod.Cells.Item(24, 3).Value = "'2"
# Customer(java.lang.Long, java.lang.String, int, int) row 25 : 0 -> 6
$wsMethod.Cells.Item(25, 3).Value = "'6"

# KafkaContainerDevMode()       row 27 : 0 -> 1
$wsMethod.Cells.Item(27, 3).Value = "'1"
